$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "67.530.30"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "3.531.20"
$ws.Range("E3").Value = "  -3.86%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell "D5" "612.29"
$ws.Range("E5").Value = "  -5.67%  "
Set-TextCell "D6" "154.30"
$ws.Range("E6").Value = "  -3.37%  "
$ws.Range("D7").Value = "3.528.02"
$ws.Range("E7").Value = "  -3.86%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextCell "D9" "0.486"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("E10").Value = "  -2.70%  "
Set-TextCell "D11" "6.90"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  -3.89%  "
Set-TextCell "D14" "32.33"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "4.124.68"
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").Value = "3.534.58"
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D17").Value = "67.488.85"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D19" "6.38"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D20" "15.59"
$ws.Range("E20").Value = "  -2.49%  "
Set-TextCell "D21" "455.02"
$ws.Range("E21").Value = "  -2.10%  "
Set-TextCell "D22" "9.41"
$ws.Range("E22").Value = "  -3.59%  "
Set-TextCell "D23" "0.643"
$ws.Range("E23").Value = "  -0.04%  "
Set-TextCell "D24" "78.72"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "3.665.86"
$ws.Range("E26").Value = "  -3.93%  "
Set-TextCell "D27" "0.0000120"
$ws.Range("E27").Value = "  -4.32%  "
Set-TextCell "D28" "10.53"
$ws.Range("E28").Value = "  -2.28%  "
Set-TextCell "D29" "8.38"
$ws.Range("E29").Value = "  -6.33%  "
Set-TextCell "D30" "1.71"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("E32").Value = "  +0.01%  "
Set-TextCell "D33" "26.05"
$ws.Range("E33").Value = "  -2.05%  "
Set-TextCell "D34" "1.92"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -3.97%  "
$ws.Range("D37").Value = "3.524.69"
$ws.Range("E37").Value = "  -3.76%  "
Set-TextCell "D38" "8.02"
$ws.Range("E38").Value = "  -4.27%  "
$ws.Range("E39").Value = "  -0.02%  "
Set-TextCell "D40" "0.999"
$ws.Range("E40").Value = "  -0.03%  "
Set-TextCell "D41" "173.47"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -4.54%  "
Set-TextCell "D43" "0.0881"
$ws.Range("E43").Value = "  -1.44%  "
Set-TextCell "D44" "2.12"
$ws.Range("E44").Value = "  -2.61%  "
Set-TextCell "D45" "0.891"
$ws.Range("E45").Value = "  -3.82%  "
Set-TextCell "D46" "29.38"
$ws.Range("E46").Value = "  +9.05%  "
Set-TextCell "D47" "45.78"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  -3.94%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell "D49" "1.23"
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D50" "7.67"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").Value = "  -2.73%  "
